$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.366.84"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.879.82"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7129"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.79"
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08043"
$ws.Range("E8").Value = "  +2.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3178"
$ws.Range("E9").Value = "  +1.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.14"
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08347"
$ws.Range("E11").Value = "  -1.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.895.81"
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.267"
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.83"
$ws.Range("E14").Value = "  +4.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7187"
$ws.Range("E15").Value = "  +0.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.372"
$ws.Range("E16").Value = "  +5.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008645"
$ws.Range("E17").Value = "  +5.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.377.80"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.34"
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.149.74"
$ws.Range("E20").Value = "  +1.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.34"
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.834"
$ws.Range("E23").Value = "  +0.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1576"
$ws.Range("E25").Value = "  -1.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.102"
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "163.33"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.63"
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.507"
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.440"
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.356"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.207"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05422"
$ws.Range("E33").Value = "  +2.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.950"
$ws.Range("E34").Value = "  +0.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7733"
$ws.Range("E35").Value = "  +3.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.189"
$ws.Range("E36").Value = "  +0.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.690"
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01893"
$ws.Range("E38").Value = "  +1.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.270.99"
$ws.Range("E39").Value = "  +3.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.752"
$ws.Range("E40").Value = "  +0.89%  "
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "113.65"
$ws.Range("E42").Value = "  +2.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9106"
$ws.Range("E43").Value = "  +2.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "74.68"
$ws.Range("E44").Value = "  +2.47%  "
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("E46").Value = "  +6.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.037.72"
$ws.Range("E47").Value = "  +0.88%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5225"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.568"
$ws.Range("E50").Value = "  +1.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4388"
$ws.Range("E51").Value = "  +1.46%  "
